$d = $word.ActiveDocument

# Locate the "git branch -av" / "- displays what branches exist..." list
# paragraph (the first occurrence in the document - the "Common Commands"
# recap list that sits right before the "1) Initializing a Local
# Repository" heading, and which uses list numId 5).
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "displays what branches exist in your repository both locally and remotely") {
        $anchorIndex = $i
        break
    }
}

$anchor = $d.Paragraphs.Item($anchorIndex)

# Duplicate the anchor paragraph (copy/paste) so the brand new paragraph
# inherits identical paragraph formatting (list numbering, spacing,
# indentation, justification) and run formatting (fonts/size/bold) without
# having to hand-build it from scratch.
$fullRange = $d.Range($anchor.Range.Start, $anchor.Range.End)
$fullRange.Copy()
$insertionPoint = $d.Range($anchor.Range.End, $anchor.Range.End)
$insertionPoint.Paste()

# The pasted paragraph is now the paragraph immediately after the anchor.
$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newRange = $newPara.Range

# Update the bold "command" run text.
$newRange.Find.Execute(
    "git branch -av ", $true, $false, $false, $false, $false, $true, 1, $false,
    "git branch -m [branch to be renamed] [new name for branch] ", 2) | Out-Null

# Update the plain-text "description" run text.
$newRange2 = $newPara.Range
$newRange2.Find.Execute(
    "- displays what branches exist in your repository both locally and remotely",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- renames a specifc branch", 2) | Out-Null
